$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"8.50811204600177e-07"
$ws.Range("B2").Value = [double]"8.839213478495367e-06"
$ws.Range("C2").Value = [double]"2.294643309141975e-06"
$ws.Range("D2").Value = [double]"1.340384824288776e-05"
$ws.Range("E2").Value = [double]"4.681493919633795e-06"
$ws.Range("F2").Value = [double]"2.035420038737357e-06"
$ws.Range("G2").Value = [double]"2.79270807368448e-06"
$ws.Range("H2").Value = [double]"2.649369662321988e-06"
$ws.Range("I2").Value = [double]"2.176171392420656e-06"
$ws.Range("J2").Value = [double]"2.083198069158243e-06"
$ws.Range("K2").Value = [double]"1.052092466125032e-05"
$ws.Range("L2").Value = [double]"3.722294650287949e-07"
$ws.Range("M2").Value = [double]"1.441560198145453e-05"
$ws.Range("N2").Value = [double]"3.44640602634172e-06"
$ws.Range("O2").Value = [double]"7.210991952888435e-06"
$ws.Range("P2").Value = [double]"2.876204234780744e-06"
$ws.Range("Q2").Value = [double]"4.467566213861573e-06"
$ws.Range("R2").Value = [double]"3.744614332390483e-06"
$ws.Range("S2").Value = [double]"2.849768634405336e-06"
$ws.Range("T2").Value = [double]"3.199328602931928e-06"
$ws.Range("U2").Value = [double]"4.420900552304374e-07"
$ws.Range("V2").Value = [double]"9.620059699955164e-07"
$ws.Range("W2").Value = [double]"4.852545771427685e-06"
$ws.Range("X2").Value = [double]"1.533534373265866e-06"
$ws.Range("Y2").Value = [double]"2.072360530291917e-06"
$ws.Range("Z2").Value = [double]"2.803835968734347e-06"
$ws.Range("AA2").Value = [double]"1.509166850155452e-06"
$ws.Range("AB2").Value = [double]"5.022390723752324e-06"
$ws.Range("AC2").Value = [double]"1.939214371304843e-06"
$ws.Range("AD2").Value = [double]"1.164019522548188e-06"
$ws.Range("AE2").Value = [double]"3.96142286263057e-06"
$ws.Range("AF2").Value = [double]"3.209864189557265e-06"
$ws.Range("AG2").Value = [double]"1.599793336026778e-06"
$ws.Range("AH2").Value = [double]"2.145271992048947e-06"
$ws.Range("AI2").Value = [double]"1.120667434406641e-06"
$ws.Range("AJ2").Value = [double]"2.063703732346767e-06"
$ws.Range("AK2").Value = [double]"4.042679393023718e-06"
$ws.Range("AL2").Value = [double]"6.070541985536693e-06"
$ws.Range("AM2").Value = [double]"5.048888169767451e-07"
$ws.Range("AN2").Value = [double]"3.3871071991598e-06"
$ws.Range("AO2").Value = [double]"1.886908989945368e-06"
$ws.Range("AP2").Value = [double]"5.354105496735428e-07"
$ws.Range("AQ2").Value = [double]"1.766446757756057e-07"
$ws.Range("AR2").Value = [double]"3.255990350226057e-06"
$ws.Range("AS2").Value = [double]"9.820552122619119e-07"
$ws.Range("AT2").Value = [double]"4.653377800423186e-06"
$ws.Range("AU2").Value = [double]"3.030621655852883e-06"
$ws.Range("AV2").Value = [double]"2.410806416719424e-07"
$ws.Range("AW2").Value = [double]"2.702173787838547e-06"
$ws.Range("AX2").Value = [double]"7.480871886400564e-07"
$ws.Range("AY2").Value = [double]"3.165809630445438e-06"
$ws.Range("AZ2").Value = [double]"2.640694901856477e-06"
$ws.Range("BA2").Value = [double]"7.631845164723927e-07"
$ws.Range("BB2").Value = [double]"1.582443360348407e-06"
$ws.Range("BC2").Value = [double]"6.252719799704209e-07"
$ws.Range("BD2").Value = [double]"7.32780245016329e-06"
$ws.Range("BE2").Value = [double]"1.372285396428197e-06"
$ws.Range("BF2").Value = [double]"7.053090484987479e-06"
$ws.Range("BG2").Value = [double]"8.730536364964792e-07"
$ws.Range("BH2").Value = [double]"1.369748702018114e-06"
$ws.Range("BI2").Value = [double]"1.813371568459843e-06"
$ws.Range("BJ2").Value = [double]"3.774544552470616e-07"
$ws.Range("BK2").Value = [double]"2.223554929514648e-06"
$ws.Range("BL2").Value = [double]"5.732806130254176e-06"
$ws.Range("BM2").Value = [double]"5.57831617697957e-06"
$ws.Range("BN2").Value = [double]"1.903240445244592e-07"
$ws.Range("BO2").Value = [double]"1.50061043768801e-06"
$ws.Range("BP2").Value = [double]"2.181400304834824e-06"
$ws.Range("BQ2").Value = [double]"2.397914329321793e-07"
$ws.Range("BR2").Value = [double]"1.634305135667091e-07"
$ws.Range("BS2").Value = [double]"1.851814772635407e-07"
$ws.Range("BT2").Value = [double]"1.30565990730247e-06"
$ws.Range("BU2").Value = [double]"3.982534508395474e-06"
$ws.Range("BV2").Value = [double]"2.431743723718682e-06"
$ws.Range("BW2").Value = [double]"6.143991413409822e-07"
$ws.Range("BX2").Value = [double]"7.185419690358685e-07"
$ws.Range("BY2").Value = [double]"5.451582637761021e-06"
$ws.Range("BZ2").Value = [double]"8.790456718088535e-08"
$ws.Range("CA2").Value = [double]"2.455874437146122e-06"
$ws.Range("CB2").Value = [double]"3.608848260228115e-07"
$ws.Range("CC2").Value = [double]"3.324342060295749e-06"
$ws.Range("CD2").Value = [double]"3.239833404222736e-06"
$ws.Range("CE2").Value = [double]"2.694625891308533e-07"
$ws.Range("CF2").Value = [double]"5.345118552213535e-07"
$ws.Range("CG2").Value = [double]"3.421109113332932e-06"
$ws.Range("CH2").Value = [double]"1.719158944979426e-06"
$ws.Range("CI2").Value = [double]"5.062010473011469e-07"
$ws.Range("CJ2").Value = [double]"4.236734412188525e-07"
$ws.Range("CK2").Value = [double]"1.02438229987456e-06"
$ws.Range("CL2").Value = [double]"1.259835471500992e-07"
$ws.Range("CM2").Value = [double]"2.684497985683265e-06"
$ws.Range("CN2").Value = [double]"2.067400146188447e-06"
$ws.Range("CO2").Value = [double]"5.219438321546477e-07"
$ws.Range("CP2").Value = [double]"1.323587639490142e-06"
$ws.Range("CQ2").Value = [double]"3.104601660197659e-07"
$ws.Range("CR2").Value = [double]"6.072277187740838e-07"
$ws.Range("CS2").Value = [double]"4.334252935223049e-06"
$ws.Range("CT2").Value = [double]"5.274448085401673e-06"
$ws.Range("CU2").Value = [double]"7.161659596022218e-06"
$ws.Range("CV2").Value = [double]"3.21104653266957e-06"
$ws.Range("CW2").Value = [double]"3.883306192165037e-07"
$ws.Range("CX2").Value = [double]"4.009680480976385e-07"
$ws.Range("CY2").Value = [double]"2.671268475751276e-06"
$ws.Range("CZ2").Value = [double]"7.591654025418393e-07"
$ws.Range("DA2").Value = [double]"1.116973749049066e-06"
$ws.Range("DB2").Value = [double]"1.126281858887523e-07"
$ws.Range("DC2").Value = [double]"1.374333010062401e-06"
$ws.Range("DD2").Value = [double]"6.442469668854756e-09"
$ws.Range("DE2").Value = [double]"8.71659494805499e-07"
$ws.Range("DF2").Value = [double]"7.498016657336848e-06"
$ws.Range("DG2").Value = [double]"4.412045655044494e-06"
$ws.Range("DH2").Value = [double]"3.247218273827457e-06"
$ws.Range("DI2").Value = [double]"4.295642611396033e-06"
$ws.Range("DJ2").Value = [double]"4.829075351153733e-06"
$ws.Range("DK2").Value = [double]"6.446736279031029e-06"
$ws.Range("DL2").Value = [double]"1.220905573973141e-06"
$ws.Range("DM2").Value = [double]"3.372478431629133e-06"
$ws.Range("DN2").Value = [double]"2.323112084923196e-06"
$ws.Range("DO2").Value = [double]"2.384856998105533e-06"
$ws.Range("DP2").Value = [double]"5.12728001922369e-06"
$ws.Range("DQ2").Value = [double]"6.521312116092304e-06"
$ws.Range("DR2").Value = [double]"3.555922376108356e-06"
$ws.Range("DS2").Value = [double]"2.101086920447415e-06"
$ws.Range("DT2").Value = [double]"6.424477874134027e-07"
$ws.Range("DU2").Value = [double]"4.94455548505357e-07"
$ws.Range("DV2").Value = [double]"1.30972409806418e-06"
$ws.Range("DW2").Value = [double]"5.233857336861547e-06"
$ws.Range("DX2").Value = [double]"2.776810788418516e-06"
$ws.Range("DY2").Value = [double]"1.138774109676888e-06"
$ws.Range("DZ2").Value = [double]"3.617408538048039e-06"
$ws.Range("EA2").Value = [double]"9.60035663410963e-07"
$ws.Range("EB2").Value = [double]"5.788498128822539e-06"
$ws.Range("EC2").Value = [double]"1.0078301215799e-07"
$ws.Range("ED2").Value = [double]"2.351581997572794e-06"
$ws.Range("EE2").Value = [double]"1.886764266600949e-06"
$ws.Range("EF2").Value = [double]"1.165568392025307e-06"
$ws.Range("EG2").Value = [double]"5.980939477012726e-07"
$ws.Range("EH2").Value = [double]"7.936753263493301e-07"
$ws.Range("EI2").Value = [double]"6.927184585947543e-06"
$ws.Range("EJ2").Value = [double]"5.416148241010887e-08"
$ws.Range("EK2").Value = [double]"6.142720394564094e-07"
$ws.Range("EL2").Value = [double]"7.24171741239843e-07"
$ws.Range("EM2").Value = [double]"2.345117763979943e-06"
$ws.Range("EN2").Value = [double]"1.242997541339719e-06"
$ws.Range("EO2").Value = [double]"4.124041879549623e-06"
$ws.Range("EP2").Value = [double]"1.642652364353125e-06"
$ws.Range("EQ2").Value = [double]"4.247401193424594e-06"
$ws.Range("ER2").Value = [double]"2.946404720205464e-06"
$ws.Range("ES2").Value = [double]"4.209217649986385e-07"
$ws.Range("ET2").Value = [double]"3.21592756336031e-06"
$ws.Range("EU2").Value = [double]"6.702102837152779e-07"
$ws.Range("EV2").Value = [double]"1.994938884308795e-06"
$ws.Range("EW2").Value = [double]"9.344657314613869e-07"
$ws.Range("EX2").Value = [double]"3.954515705117956e-06"
$ws.Range("EY2").Value = [double]"4.074821049471211e-07"
$ws.Range("EZ2").Value = [double]"4.565748099594202e-07"
$ws.Range("FA2").Value = [double]"3.253088380006375e-06"
$ws.Range("FB2").Value = [double]"1.589070052432362e-06"
$ws.Range("FC2").Value = [double]"4.699492706095043e-07"
$ws.Range("FD2").Value = [double]"5.083195446786704e-07"
$ws.Range("FE2").Value = [double]"3.791067513247981e-07"
$ws.Range("FF2").Value = [double]"7.380143074442458e-07"
$ws.Range("FG2").Value = [double]"1.503799239799264e-06"
$ws.Range("FH2").Value = [double]"4.631281171896262e-06"
$ws.Range("FI2").Value = [double]"6.252957973629236e-07"
$ws.Range("FJ2").Value = [double]"3.649385234893998e-06"
$ws.Range("FK2").Value = [double]"2.608967406558804e-06"
$ws.Range("FL2").Value = [double]"2.224645413662074e-06"
$ws.Range("FM2").Value = [double]"2.001529992412543e-06"
$ws.Range("FN2").Value = [double]"8.100167292468541e-07"
$ws.Range("FO2").Value = [double]"1.580579237270285e-06"
$ws.Range("FP2").Value = [double]"7.467054274457041e-06"
$ws.Range("FQ2").Value = [double]"5.369371137931012e-06"
$ws.Range("FR2").Value = [double]"2.308399587036547e-07"
$ws.Range("FS2").Value = [double]"1.651900674914941e-06"
$ws.Range("FT2").Value = [double]"2.195393108195276e-06"
$ws.Range("FU2").Value = [double]"1.200508563670155e-06"
$ws.Range("FV2").Value = [double]"1.212611641676631e-05"
$ws.Range("FW2").Value = [double]"5.320078798831673e-06"
$ws.Range("FX2").Value = [double]"2.25535586650949e-06"
$ws.Range("FY2").Value = [double]"6.899260370119009e-06"
$ws.Range("FZ2").Value = [double]"6.698894026158086e-07"
$ws.Range("GA2").Value = [double]"7.617085884703556e-06"
$ws.Range("GB2").Value = [double]"4.231662842357764e-08"
$ws.Range("GC2").Value = [double]"4.547125172393862e-06"
$ws.Range("GD2").Value = [double]"5.624478944810107e-06"
$ws.Range("GE2").Value = [double]"5.502164185600122e-06"
$ws.Range("GF2").Value = [double]"1.048925355462416e-06"
$ws.Range("GG2").Value = [double]"2.280168018842232e-06"
